$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.870.67"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.544.83"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0829"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "2.938.74"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "2.564.41"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.872"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "42.911.12"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.64%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0798"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.37%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +27.87%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "2.091.20"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0305"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("D49").Value = "2.796.48"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +8.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "
